# binding_affinity.xlsx update
# - fixes number of significant digits so SE digit count matches the mean's
# - converts a couple of AVERAGE/STDEVA formulas to their plain computed values
# - clears a stray fill flag on the "Tessler et al., 1994" row
# - updates the saved selection / active sheet to reflect where the author ended up

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("VEGFA165_VEGFR1")
$ws2 = $wb.Worksheets.Item("VEGFA165_VEGFR2")
$ws3 = $wb.Worksheets.Item("VEGFA165_NRP1")

# --- VEGFA165_VEGFR1 -------------------------------------------------
# Whitaker et al., 2001 (COS-1 cell) SE: 0.96 -> 0.956
$ws1.Range("D7").Value = 0.95599999999999996

# --- VEGFA165_VEGFR2 ---------------------------------------------------
# Huang et al., 1998 (VEGF-A165): drop the AVERAGE/STDEVA formulas, keep
# plain numbers (rounded mean/SE pair updated together)
$ws2.Range("C2").Value = 230
$ws2.Range("D2").Value = 120

# Huang et al., 1998 (VEGF-A164): same treatment
$ws2.Range("C3").Value = 240
$ws2.Range("D3").Value = 100

# Lu et al., 2023 mean/SE significant digits
$ws2.Range("C15").Value = 115.4
$ws2.Range("D15").Value = 73.44

# Tessler et al., 1994 row: clear the leftover fill flag on A18:C18 so the
# cells fall back to the same (unfilled) style used elsewhere in the table
$ws2.Range("A18:C18").Interior.Pattern = -4142

# --- VEGFA165_NRP1 -----------------------------------------------------
# Waltenberger et al., 1994 (HUVEC) SE: 0.82 -> 0.818
$ws3.Range("D2").Value = 0.81799999999999995

# Lu et al., 2023 mean/SE significant digits
$ws3.Range("C11").Value = 0.14510000000000001
$ws3.Range("D11").Value = 0.058909999999999997

# --- Selection / active-sheet bookkeeping ------------------------------
# Set the stored selection on each sheet, activating the last one
# (VEGFA165_VEGFR2) so it ends up as the active tab on save.
$ws1.Range("D8").Select()
$ws3.Range("D12").Select()
$ws2.Range("C22").Select()
